# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Update record_atd (column C) values for the affected rows and
# refresh the dependent record_id (column D) values that mirror/average
# them, plus the column C summary average in row 78.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 18
$ws.Range("D3").Value = 18
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 22
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 38
$ws.Range("C8").Value = 105
$ws.Range("D8").Value = 105
$ws.Range("C10").Value = 225
$ws.Range("D10").Value = 225
$ws.Range("C12").Value = 64
$ws.Range("D12").Value = 64
$ws.Range("C14").Value = 26
$ws.Range("D14").Value = 26
$ws.Range("C16").Value = 152
$ws.Range("D16").Value = 152
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 10
$ws.Range("C20").Value = 53
$ws.Range("D20").Value = 53
$ws.Range("C22").Value = 198
$ws.Range("D22").Value = 198
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("C26").Value = 151
$ws.Range("D26").Value = 151
$ws.Range("C28").Value = 30
$ws.Range("D28").Value = 30
$ws.Range("C29").Value = 114
$ws.Range("D29").Value = 114
$ws.Range("C31").Value = 21
$ws.Range("D31").Value = 21
$ws.Range("C33").Value = 84
$ws.Range("D33").Value = 84
$ws.Range("C35").Value = 170
$ws.Range("D35").Value = 170
$ws.Range("C37").Value = 12
$ws.Range("D37").Value = 12
$ws.Range("C39").Value = 174
$ws.Range("D39").Value = 174
$ws.Range("C41").Value = 106
$ws.Range("D41").Value = 106
$ws.Range("C44").Value = 29
$ws.Range("D44").Value = 131
$ws.Range("C45").Value = 36
$ws.Range("D45").Value = 36
$ws.Range("C46").Value = 234
$ws.Range("D46").Value = 234
$ws.Range("C49").Value = 115
$ws.Range("D49").Value = 115
$ws.Range("C50").Value = 111
$ws.Range("D50").Value = 111
$ws.Range("C52").Value = 208
$ws.Range("D52").Value = 208
$ws.Range("C54").Value = 9
$ws.Range("D54").Value = 9
$ws.Range("C56").Value = 177
$ws.Range("D56").Value = 177
$ws.Range("C58").Value = 17
$ws.Range("D58").Value = 17
$ws.Range("C61").Value = 42
$ws.Range("D61").Value = 42
$ws.Range("C62").Value = 19
$ws.Range("D62").Value = 19
$ws.Range("C64").Value = 155
$ws.Range("D64").Value = 155
$ws.Range("C66").Value = 67
$ws.Range("D66").Value = 67
$ws.Range("C69").Value = 41
$ws.Range("D69").Value = 41
$ws.Range("C70").Value = 68
$ws.Range("D70").Value = 68
$ws.Range("C72").Value = 45
$ws.Range("D72").Value = 36
$ws.Range("C74").Value = 11
$ws.Range("D74").Value = 11
$ws.Range("C75").Value = 150
$ws.Range("D75").Value = 150
$ws.Range("C77").Value = 113
$ws.Range("D77").Value = 113
$ws.Range("C78").Value = 85.15000000000001
